# Refresh Universalis market-data columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) across all eight Mateus job sheets, per the scheduled-runner sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 221.14285
$ws.Range("I9").Value = 183
$ws.Range("K9").Value = 183
$ws.Range("M9").Value = -14
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H70").Value = 4041.3845
$ws.Range("I70").Value = 1871
$ws.Range("J70").Value = 4436
$ws.Range("K70").Value = 5613
$ws.Range("L70").Value = 13308
$ws.Range("M70").Value = -5343
$ws.Range("N70").Value = -13848
$ws.Range("H73").Value = 4041.3845
$ws.Range("I73").Value = 1871
$ws.Range("J73").Value = 4436
$ws.Range("K73").Value = 5613
$ws.Range("L73").Value = 13308
$ws.Range("M73").Value = -4677
$ws.Range("N73").Value = -15180
$ws.Range("H86").Value = 3120.6924
$ws.Range("I86").Value = 2457.3
$ws.Range("K86").Value = 2457.3
$ws.Range("M86").Value = -1334.3
$ws.Range("H89").Value = 3120.6924
$ws.Range("I89").Value = 2457.3
$ws.Range("K89").Value = 12286.5
$ws.Range("M89").Value = -6670.5
$ws.Range("H98").Value = 861.7241
$ws.Range("I98").Value = 869.03845
$ws.Range("K98").Value = 869.03845
$ws.Range("M98").Value = 628.96155
$ws.Range("H99").Value = 206.25
$ws.Range("I99").Value = 191.66667
$ws.Range("K99").Value = 575.00001
$ws.Range("M99").Value = 922.99999
$ws.Range("H103").Value = 464.14285
$ws.Range("I103").Value = 399.75
$ws.Range("K103").Value = 1199.25
$ws.Range("M103").Value = -613.25
$ws.Range("H122").Value = 861.7241
$ws.Range("I122").Value = 869.03845
$ws.Range("K122").Value = 2607.11535
$ws.Range("M122").Value = -157.11535
$ws.Range("H127").Value = 1097
$ws.Range("I127").Value = 1116.4
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 3349.2
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 1610.8
$ws.Range("N127").Value = -12920
$ws.Range("H129").Value = 1225.1
$ws.Range("I129").Value = 1281.375
$ws.Range("J129").Value = 1000
$ws.Range("K129").Value = 3844.125
$ws.Range("L129").Value = 3000
$ws.Range("M129").Value = 1155.875
$ws.Range("N129").Value = -13000
$ws.Range("H131").Value = 503000
$ws.Range("I131").Value = 503000
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1509000
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = -1503960
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 1171.5897
$ws.Range("I132").Value = 1010.44116
$ws.Range("J132").Value = 2267.4
$ws.Range("K132").Value = 3031.32348
$ws.Range("L132").Value = 6802.200000000001
$ws.Range("M132").Value = -501.32348
$ws.Range("N132").Value = -11862.2
$ws.Range("H135").Value = 560.9231
$ws.Range("I135").Value = 483.36
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 4350.24
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -1815.24
$ws.Range("N135").Value = -27570
$ws.Range("H137").Value = 1576.6
$ws.Range("I137").Value = 1446.1333
$ws.Range("J137").Value = 1968
$ws.Range("K137").Value = 4338.3999
$ws.Range("L137").Value = 5904
$ws.Range("M137").Value = -1788.3999
$ws.Range("N137").Value = -11004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2055.5454
$ws.Range("I32").Value = 2113.0952
$ws.Range("K32").Value = 2113.0952
$ws.Range("M32").Value = -1826.0952
$ws.Range("H74").Value = 3204.4048
$ws.Range("I74").Value = 2579.3572
$ws.Range("K74").Value = 2579.3572
$ws.Range("M74").Value = -1705.3572
$ws.Range("H77").Value = 3204.4048
$ws.Range("I77").Value = 2579.3572
$ws.Range("K77").Value = 12896.786
$ws.Range("M77").Value = -8528.786
$ws.Range("H110").Value = 6835.9414
$ws.Range("I110").Value = 4267.6665
$ws.Range("K110").Value = 4267.6665
$ws.Range("M110").Value = -2222.6665
$ws.Range("H122").Value = 2764.08
$ws.Range("I122").Value = 2358.1765
$ws.Range("J122").Value = 3626.625
$ws.Range("K122").Value = 7074.529500000001
$ws.Range("L122").Value = 10879.875
$ws.Range("M122").Value = -4624.529500000001
$ws.Range("N122").Value = -15779.875
$ws.Range("H134").Value = 139989.67
$ws.Range("J134").Value = 139989.67
$ws.Range("L134").Value = 139989.67
$ws.Range("N134").Value = -150129.67
$ws.Range("H139").Value = 94040.664
$ws.Range("J139").Value = 121061
$ws.Range("L139").Value = 121061
$ws.Range("N139").Value = -131341

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H107").Value = 1744.9166
$ws.Range("I107").Value = 1704.1904
$ws.Range("K107").Value = 1704.1904
$ws.Range("M107").Value = 215.8096
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3657.7646
$ws.Range("I16").Value = 2164.1428
$ws.Range("J16").Value = 4703.3
$ws.Range("K16").Value = 2164.1428
$ws.Range("L16").Value = 4703.3
$ws.Range("M16").Value = -1877.1428
$ws.Range("N16").Value = -5277.3
$ws.Range("H113").Value = 3657.7646
$ws.Range("I113").Value = 2164.1428
$ws.Range("J113").Value = 4703.3
$ws.Range("K113").Value = 2164.1428
$ws.Range("L113").Value = 4703.3
$ws.Range("M113").Value = 5.857199999999921
$ws.Range("N113").Value = -9043.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 500251
$ws.Range("I36").Value = 500251
$ws.Range("K36").Value = 1500753
$ws.Range("M36").Value = -1500584
$ws.Range("H51").Value = 34120.777
$ws.Range("I51").Value = 14999.25
$ws.Range("K51").Value = 44997.75
$ws.Range("M51").Value = -44537.75
$ws.Range("H128").Value = 917231.25
$ws.Range("I128").Value = 917231.25
$ws.Range("K128").Value = 2751693.75
$ws.Range("M128").Value = -2746713.75
$ws.Range("H138").Value = 2002
$ws.Range("I138").Value = 1669.6666
$ws.Range("K138").Value = 5008.9998
$ws.Range("M138").Value = 131.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1443.8572
$ws.Range("I113").Value = 1419.5
$ws.Range("K113").Value = 1419.5
$ws.Range("M113").Value = 750.5
$ws.Range("H122").Value = 2071.5642
$ws.Range("I122").Value = 2083.1482
$ws.Range("J122").Value = 2045.5
$ws.Range("K122").Value = 6249.444600000001
$ws.Range("L122").Value = 6136.5
$ws.Range("M122").Value = -3799.444600000001
$ws.Range("N122").Value = -11036.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1430
$ws.Range("I22").Value = 1712.5
$ws.Range("K22").Value = 1712.5
$ws.Range("M22").Value = -1417.5
$ws.Range("H27").Value = 1430
$ws.Range("I27").Value = 1712.5
$ws.Range("K27").Value = 1712.5
$ws.Range("M27").Value = -1605.5
$ws.Range("H42").Value = 13749.5
$ws.Range("J42").Value = 14999.333
$ws.Range("L42").Value = 14999.333
$ws.Range("N42").Value = -16125.333
$ws.Range("H46").Value = 8904.174000000001
$ws.Range("I46").Value = 2964.1428
$ws.Range("J46").Value = 13893.8
$ws.Range("K46").Value = 2964.1428
$ws.Range("L46").Value = 13893.8
$ws.Range("M46").Value = -2776.1428
$ws.Range("N46").Value = -14269.8
$ws.Range("H49").Value = 13749.5
$ws.Range("J49").Value = 14999.333
$ws.Range("L49").Value = 14999.333
$ws.Range("N49").Value = -15293.333
$ws.Range("H55").Value = 699.9286
$ws.Range("I55").Value = 816.6667
$ws.Range("J55").Value = 489.8
$ws.Range("K55").Value = 816.6667
$ws.Range("L55").Value = 489.8
$ws.Range("M55").Value = -643.6667
$ws.Range("N55").Value = -835.8
$ws.Range("H68").Value = 1991.75
$ws.Range("I68").Value = 1999
$ws.Range("J68").Value = 1989.3334
$ws.Range("K68").Value = 1999
$ws.Range("L68").Value = 1989.3334
$ws.Range("M68").Value = -1250
$ws.Range("N68").Value = -3487.3334
$ws.Range("H71").Value = 1991.75
$ws.Range("I71").Value = 1999
$ws.Range("J71").Value = 1989.3334
$ws.Range("K71").Value = 9995
$ws.Range("L71").Value = 9946.666999999999
$ws.Range("M71").Value = -6251
$ws.Range("N71").Value = -17434.667
$ws.Range("H122").Value = 4064.7
$ws.Range("I122").Value = 3798.2
$ws.Range("J122").Value = 4331.2
$ws.Range("K122").Value = 11394.6
$ws.Range("L122").Value = 12993.6
$ws.Range("M122").Value = -8944.599999999999
$ws.Range("N122").Value = -17893.6
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4261.25
$ws.Range("I122").Value = 4057.8235
$ws.Range("J122").Value = 4755.2856
$ws.Range("K122").Value = 12173.4705
$ws.Range("L122").Value = 14265.8568
$ws.Range("M122").Value = -9723.470499999999
$ws.Range("N122").Value = -19165.8568
$ws.Range("H132").Value = 8214.143
$ws.Range("I132").Value = 7499.8
$ws.Range("K132").Value = 22499.4
$ws.Range("M132").Value = -19969.4
